# Scheduled-runner style refresh of market-price-derived profit figures
# across the Sargatanas_Profits workbook's per-job sheets (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR). Updates currentAveragePrice*/LevePrice*/
# LeveProfit* columns (H-N) for the specific leve rows that moved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 830.3333
$ws.Range("I29").Value = 830.3333
$ws.Range("K29").Value = 2490.9999
$ws.Range("M29").Value = -2209.9999

$ws.Range("H40").Value = 3093.5715
$ws.Range("I40").Value = 3285.7144
$ws.Range("J40").Value = 2901.4285
$ws.Range("K40").Value = 3285.7144
$ws.Range("L40").Value = 2901.4285
$ws.Range("M40").Value = -3110.7144
$ws.Range("N40").Value = -3251.4285

$ws.Range("H61").Value = 608.6667
$ws.Range("I61").Value = 608.6667
$ws.Range("K61").Value = 1826.0001
$ws.Range("M61").Value = -1654.0001

$ws.Range("H74").Value = 71435816
$ws.Range("I74").Value = 83338450
$ws.Range("K74").Value = 83338450
$ws.Range("M74").Value = -83337514

$ws.Range("H77").Value = 71435816
$ws.Range("I77").Value = 83338450
$ws.Range("K77").Value = 416692250
$ws.Range("M77").Value = -416687570

$ws.Range("H87").Value = 60000
$ws.Range("J87").Value = 60000
$ws.Range("L87").Value = 60000
$ws.Range("N87").Value = -62496

$ws.Range("H90").Value = 60000
$ws.Range("J90").Value = 60000
$ws.Range("L90").Value = 180000
$ws.Range("N90").Value = -192480

$ws.Range("H100").Value = 2665.6365
$ws.Range("I100").Value = 2144.6667
$ws.Range("K100").Value = 2144.6667
$ws.Range("M100").Value = -1603.6667

$ws.Range("H113").Value = 35953856
$ws.Range("J113").Value = 45462764
$ws.Range("L113").Value = 45462764
$ws.Range("N113").Value = -45469272

$ws.Range("H116").Value = 31261186
$ws.Range("I116").Value = 83338170
$ws.Range("K116").Value = 83338170
$ws.Range("M116").Value = -83334728

$ws.Range("H132").Value = 2037.3572
$ws.Range("I132").Value = 1647.9584
$ws.Range("J132").Value = 4373.75
$ws.Range("K132").Value = 4943.8752
$ws.Range("L132").Value = 13121.25
$ws.Range("M132").Value = -2413.8752
$ws.Range("N132").Value = -18181.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1738408
$ws.Range("I32").Value = 1787791
$ws.Range("K32").Value = 1787791
$ws.Range("M32").Value = -1787504

$ws.Range("H102").Value = 1169.1875
$ws.Range("J102").Value = 1283.3334
$ws.Range("L102").Value = 1283.3334
$ws.Range("N102").Value = -4527.3334

$ws.Range("H132").Value = 6312.6978
$ws.Range("I132").Value = 5048.3667
$ws.Range("K132").Value = 15145.1001
$ws.Range("M132").Value = -12615.1001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2620.9565
$ws.Range("J94").Value = 4910.5557
$ws.Range("L94").Value = 4910.5557
$ws.Range("N94").Value = -5812.5557

$ws.Range("H96").Value = 30453.5
$ws.Range("I96").Value = 16770.2
$ws.Range("K96").Value = 16770.2
$ws.Range("M96").Value = -14024.2

$ws.Range("H134").Value = 6650.4053
$ws.Range("I134").Value = 2915.2354
$ws.Range("K134").Value = 8745.706200000001
$ws.Range("M134").Value = -6210.706200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8357.171
$ws.Range("J31").Value = 12088.8
$ws.Range("L31").Value = 12088.8
$ws.Range("N31").Value = -12678.8

$ws.Range("H34").Value = 8357.171
$ws.Range("J34").Value = 12088.8
$ws.Range("L34").Value = 12088.8
$ws.Range("N34").Value = -12492.8

$ws.Range("H37").Value = 11057
$ws.Range("J37").Value = 11057
$ws.Range("L37").Value = 11057
$ws.Range("N37").Value = -11271

$ws.Range("H58").Value = 10006068
$ws.Range("I58").Value = 22730206
$ws.Range("J58").Value = 8529.429
$ws.Range("K58").Value = 22730206
$ws.Range("L58").Value = 8529.429
$ws.Range("M58").Value = -22730003
$ws.Range("N58").Value = -8935.429

$ws.Range("H134").Value = 8234
$ws.Range("I134").Value = 4403.9165
$ws.Range("K134").Value = 13211.7495
$ws.Range("M134").Value = -10676.7495

$ws.Range("H136").Value = 10006068
$ws.Range("I136").Value = 22730206
$ws.Range("J136").Value = 8529.429
$ws.Range("K136").Value = 68190618
$ws.Range("L136").Value = 25588.287
$ws.Range("M136").Value = -68188068
$ws.Range("N136").Value = -30688.287

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 1194.9
$ws.Range("I15").Value = 10.5
$ws.Range("J15").Value = 1984.5
$ws.Range("K15").Value = 31.5
$ws.Range("L15").Value = 5953.5
$ws.Range("M15").Value = 108.5
$ws.Range("N15").Value = -6233.5

$ws.Range("H18").Value = 1833.3334
$ws.Range("I18").Value = 1833.3334
$ws.Range("K18").Value = 5500.0002
$ws.Range("M18").Value = -5331.0002

$ws.Range("H34").Value = 5410.8335
$ws.Range("I34").Value = 209
$ws.Range("J34").Value = 6897.0713
$ws.Range("K34").Value = 627
$ws.Range("L34").Value = 20691.2139
$ws.Range("M34").Value = -543
$ws.Range("N34").Value = -20859.2139

$ws.Range("H59").Value = 999
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H64").Value = 20002540
$ws.Range("I64").Value = 2566.3333
$ws.Range("K64").Value = 7698.999899999999
$ws.Range("M64").Value = -7428.999899999999

$ws.Range("H67").Value = 20002540
$ws.Range("I67").Value = 2566.3333
$ws.Range("K67").Value = 7698.999899999999
$ws.Range("M67").Value = -6762.999899999999

$ws.Range("H122").Value = 4098.6665
$ws.Range("I122").Value = 2177.8
$ws.Range("J122").Value = 6499.75
$ws.Range("K122").Value = 19600.2
$ws.Range("L122").Value = 58497.75
$ws.Range("M122").Value = -17150.2
$ws.Range("N122").Value = -63397.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 109663
$ws.Range("J134").Value = 109663
$ws.Range("L134").Value = 328989
$ws.Range("N134").Value = -334059

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7351.8237
$ws.Range("I7").Value = 6530.3335
$ws.Range("J7").Value = 8276
$ws.Range("K7").Value = 6530.3335
$ws.Range("L7").Value = 8276
$ws.Range("M7").Value = -6418.3335
$ws.Range("N7").Value = -8500

$ws.Range("H46").Value = 3258.3635
$ws.Range("I46").Value = 2450
$ws.Range("J46").Value = 3720.2856
$ws.Range("K46").Value = 2450
$ws.Range("L46").Value = 3720.2856
$ws.Range("M46").Value = -2262
$ws.Range("N46").Value = -4096.2856

$ws.Range("H93").Value = 5978.3335
$ws.Range("I93").Value = 5706.25
$ws.Range("J93").Value = 7066.6665
$ws.Range("K93").Value = 5706.25
$ws.Range("L93").Value = 7066.6665
$ws.Range("M93").Value = -4458.25
$ws.Range("N93").Value = -9562.666499999999

$ws.Range("H126").Value = 7351.8237
$ws.Range("I126").Value = 6530.3335
$ws.Range("J126").Value = 8276
$ws.Range("K126").Value = 19591.0005
$ws.Range("L126").Value = 24828
$ws.Range("M126").Value = -17121.0005
$ws.Range("N126").Value = -29768

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 37079730
$ws.Range("I136").Value = 125001416
$ws.Range("K136").Value = 375004248
$ws.Range("M136").Value = -375001698
